$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header cells J2 ("сумм") / K2 ("тк") - styled like the other header cells
# (bold font, center/center, wrap) but with only a thick LEFT border instead
# of a full thick box border.
# ---------------------------------------------------------------------------
$ws.Range("J2").Value = "сумм"
$ws.Range("K2").Value = "тк"

$ws.Range("J2:K2").Font.Bold = $true
$ws.Range("J2:K2").Font.Size = 10
$ws.Range("J2:K2").Font.Name = "Arial"
$ws.Range("J2:K2").Font.Color = $ws.Range("C2").Font.Color
$ws.Range("J2:K2").HorizontalAlignment = -4108
$ws.Range("J2:K2").VerticalAlignment = -4108
$ws.Range("J2:K2").WrapText = $true
$ws.Range("J2:K2").Borders.Item(7).Weight = 4

# ---------------------------------------------------------------------------
# J4:J31 - row total (homework points), K4:K31 - "тк" grade derived from it.
# J4 is entered on its own; J5:J31 filled as one block so Excel stores it as
# a shared formula group (matches how the sheet was actually authored).
# ---------------------------------------------------------------------------
$ws.Range("J4").Formula = "=SUM(C4:H4)"
$ws.Range("J5:J31").Formula = "=SUM(C5:H5)"

$kValues = @{
    4 = 3;  5 = 5;  6 = 4;  7 = 3;  8 = 4;  9 = 5;  10 = 5;
    11 = 3; 12 = 5; 13 = 5; 14 = 3; 15 = 5; 16 = 5; 17 = 5;
    18 = 5; 19 = 3; 20 = 5; 21 = 3; 22 = 3; 23 = 5; 24 = 5;
    25 = 4; 26 = 4; 27 = 3; 28 = 3; 29 = 5; 30 = 5; 31 = 5;
}
foreach ($row in 4..31) {
    $ws.Cells.Item($row, 11).Value = $kValues[$row]
}

# ---------------------------------------------------------------------------
# Conditional formatting (3-colour scale) over the new "сумм" column.
# ---------------------------------------------------------------------------
$cf = $ws.Range("J4:J31").FormatConditions.AddColorScale(3)
$cf.ColorScaleCriteria.Item(1).Type = 4        # xlConditionValueLowestValue
$cf.ColorScaleCriteria.Item(1).FormatColor.Color = 7039080     # FFF8696B (BGR)
$cf.ColorScaleCriteria.Item(2).Type = 3        # xlConditionValuePercentile
$cf.ColorScaleCriteria.Item(2).Value = 50
$cf.ColorScaleCriteria.Item(2).FormatColor.Color = 8684223     # FFFFEB84 (BGR)
$cf.ColorScaleCriteria.Item(3).Type = 5        # xlConditionValueHighestValue
$cf.ColorScaleCriteria.Item(3).FormatColor.Color = 8109411     # FF63BE7B (BGR)

# ---------------------------------------------------------------------------
# Selection, matching the author's last recorded cursor position.
# ---------------------------------------------------------------------------
$ws.Range("K4").Select()

Write-Output "edit complete"
